$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$sh1 = $np.Shapes.Item(1)
$sh2 = $np.Shapes.Item(2)
$sh3 = $np.Shapes.Item(3)
$sh2.TextFrame.TextRange.Text = "Hello notes"
Write-Host "Shapes on notes page after write to shape2, having touched 1 and 3 first:"
for ($i = 1; $i -le $np.Shapes.Count; $i++) {
    $sh = $np.Shapes.Item($i)
    Write-Host "  $i : $($sh.Name)"
}
